$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "IsActive" header in cell C1, matching existing headers A1/B1
$ws.Range("C1").Value = "IsActive"
